$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:AF2").ClearContents()
$ws.Range("AI2").ClearContents()
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 1.23
$ws.Range("AJ2").Value = 93020173

# Row 3
$ws.Range("D3:AF3").ClearContents()
$ws.Range("AI3").ClearContents()
$ws.Range("AG3").Value = 300
$ws.Range("AH3").Value = 1.71
$ws.Range("AJ3").Value = 93020173

# Row 4
$ws.Range("Y4:Z4").ClearContents()
$ws.Range("D4").Value = 6900
$ws.Range("E4").Value = 2616
$ws.Range("F4").Value = 2616
$ws.Range("G4").Value = 3019
$ws.Range("H4").Value = 2695
$ws.Range("I4").Value = 2615
$ws.Range("J4").Value = 79
$ws.Range("K4").Value = 30975
$ws.Range("L4").Value = 2281
$ws.Range("M4").Value = 28694
$ws.Range("N4").Value = 28034
$ws.Range("O4").Value = 659
$ws.Range("P4").Value = 465
$ws.Range("Q4").Value = 858
$ws.Range("R4").Value = 749
$ws.Range("S4").Value = -605
$ws.Range("T4").Value = 110
$ws.Range("U4").Value = 749
$ws.Range("V4").Value = 295
$ws.Range("W4").Value = 37.91
$ws.Range("X4").Value = 39.05
$ws.Range("AA4").Value = 7.95
$ws.Range("AB4").Value = 5991.64
$ws.Range("AC4").Value = 2812
$ws.Range("AD4").Value = 7.38
$ws.Range("AE4").Value = 30573
$ws.Range("AF4").Value = 0.68
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 1.45
$ws.Range("AI4").Value = 10.52
$ws.Range("AJ4").Value = 93020173

# Row 5
$ws.Range("D5").Value = 8248
$ws.Range("E5").Value = 1993
$ws.Range("F5").Value = 1993
$ws.Range("G5").Value = 2727
$ws.Range("H5").Value = 2425
$ws.Range("I5").Value = 2331
$ws.Range("J5").Value = 94
$ws.Range("K5").Value = 33161
$ws.Range("L5").Value = 2564
$ws.Range("M5").Value = 30597
$ws.Range("N5").Value = 29846
$ws.Range("O5").Value = 750
$ws.Range("P5").Value = 465
$ws.Range("Q5").Value = 1229
$ws.Range("R5").Value = -2129
$ws.Range("S5").Value = -60
$ws.Range("T5").Value = 232
$ws.Range("U5").Value = 997
$ws.Range("V5").Value = 500
$ws.Range("W5").Value = 24.16
$ws.Range("X5").Value = 29.4
$ws.Range("Y5").Value = 8.050000000000001
$ws.Range("Z5").Value = 7.56
$ws.Range("AA5").Value = 8.380000000000001
$ws.Range("AB5").Value = 6449.13
$ws.Range("AC5").Value = 2506
$ws.Range("AD5").Value = 7.54
$ws.Range("AE5").Value = 32549
$ws.Range("AF5").Value = 0.58
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 1.59
$ws.Range("AI5").Value = 11.8
$ws.Range("AJ5").Value = 93020173

# Row 6
$ws.Range("D6").Value = 8480
$ws.Range("E6").Value = 2115
$ws.Range("F6").Value = 2115
$ws.Range("G6").Value = 2216
$ws.Range("H6").Value = 2184
$ws.Range("I6").Value = 2071
$ws.Range("K6").Value = 34986
$ws.Range("L6").Value = 2627
$ws.Range("M6").Value = 32359
$ws.Range("N6").Value = 31493
$ws.Range("P6").Value = 465
$ws.Range("Q6").Value = 609
$ws.Range("R6").Value = -17
$ws.Range("S6").Value = 96
$ws.Range("T6").Value = 213
$ws.Range("U6").Value = 396
$ws.Range("V6").Value = 864
$ws.Range("W6").Value = 24.94
$ws.Range("X6").Value = 25.75
$ws.Range("Y6").Value = 6.75
$ws.Range("Z6").Value = 6.41
$ws.Range("AA6").Value = 8.119999999999999
$ws.Range("AB6").Value = 6829.65
$ws.Range("AC6").Value = 2226
$ws.Range("AD6").Value = 7.46
$ws.Range("AE6").Value = 34346
$ws.Range("AF6").Value = 0.48
$ws.Range("AG6").Value = 300
$ws.Range("AH6").Value = 1.81
$ws.Range("AI6").Value = 13.28
$ws.Range("AJ6").Value = 93020173

# Row 7
$ws.Range("D7").Value = 8779
$ws.Range("E7").Value = 2203
$ws.Range("G7").Value = 2435
$ws.Range("H7").Value = 2100
$ws.Range("I7").Value = 1992
$ws.Range("K7").Value = 36935
$ws.Range("L7").Value = 2660
$ws.Range("M7").Value = 34278
$ws.Range("N7").Value = 33330
$ws.Range("P7").Value = 467
$ws.Range("Q7").Value = 1819
$ws.Range("R7").Value = -934
$ws.Range("S7").Value = -270
$ws.Range("T7").Value = 220
$ws.Range("U7").Value = 1614
$ws.Range("W7").Value = 25.09
$ws.Range("X7").Value = 23.93
$ws.Range("Y7").Value = 6.15
$ws.Range("Z7").Value = 5.84
$ws.Range("AA7").Value = 7.76
$ws.Range("AC7").Value = 2142
$ws.Range("AD7").Value = 6.19
$ws.Range("AE7").Value = 36348
$ws.Range("AF7").Value = 0.36
$ws.Range("AG7").Value = 300
$ws.Range("AH7").Value = 2.26
$ws.Range("AI7").Value = 14.01

# Row 8
$ws.Range("D8").Value = 8995
$ws.Range("E8").Value = 2406
$ws.Range("G8").Value = 2537
$ws.Range("H8").Value = 2223
$ws.Range("I8").Value = 2098
$ws.Range("K8").Value = 38886
$ws.Range("L8").Value = 2595
$ws.Range("M8").Value = 36291
$ws.Range("N8").Value = 35252
$ws.Range("P8").Value = 467
$ws.Range("Q8").Value = 1985
$ws.Range("R8").Value = -888
$ws.Range("S8").Value = -266
$ws.Range("T8").Value = 258
$ws.Range("U8").Value = 1307
$ws.Range("W8").Value = 26.74
$ws.Range("X8").Value = 24.72
$ws.Range("Y8").Value = 6.12
$ws.Range("Z8").Value = 5.87
$ws.Range("AA8").Value = 7.15
$ws.Range("AC8").Value = 2256
$ws.Range("AD8").Value = 5.87
$ws.Range("AE8").Value = 38444
$ws.Range("AF8").Value = 0.34
$ws.Range("AG8").Value = 300
$ws.Range("AH8").Value = 2.26
$ws.Range("AI8").Value = 13.3

# Row 9
$ws.Range("D9").Value = 9246
$ws.Range("E9").Value = 2472
$ws.Range("G9").Value = 2614
$ws.Range("H9").Value = 2286
$ws.Range("I9").Value = 2158
$ws.Range("K9").Value = 40881
$ws.Range("L9").Value = 2504
$ws.Range("M9").Value = 38376
$ws.Range("N9").Value = 37243
$ws.Range("P9").Value = 467
$ws.Range("Q9").Value = 2032
$ws.Range("R9").Value = -778
$ws.Range("S9").Value = -264
$ws.Range("T9").Value = 265
$ws.Range("U9").Value = 1429
$ws.Range("W9").Value = 26.74
$ws.Range("X9").Value = 24.73
$ws.Range("Y9").Value = 5.95
$ws.Range("Z9").Value = 5.73
$ws.Range("AA9").Value = 6.53
$ws.Range("AC9").Value = 2320
$ws.Range("AD9").Value = 5.71
$ws.Range("AE9").Value = 40616
$ws.Range("AF9").Value = 0.33
$ws.Range("AG9").Value = 300
$ws.Range("AH9").Value = 2.26
$ws.Range("AI9").Value = 12.93
